# Clean unused freezing demo variables
# Removes the "freeze_demo" / "time_demo_freezing" rows (rows 47 and 48) from
# the "constant" worksheet, shifting everything below them up by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")
$ws.Activate()

# Delete the two rows that held the now-unused "freeze_demo" and
# "time_demo_freezing" settings; this shifts rows 49:57 up to 47:55 and
# automatically adjusts the data validation ranges that reference them.
$ws.Rows("47:48").Delete()

# Reflect the new selection left behind by the edit.
$ws.Range("A52").Select()
